# Auto-generated edit script for Sagittarius_Profits workbook update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 505.7143
$ws.Range("J58").Value = 610
$ws.Range("L58").Value = 1830
$ws.Range("N58").Value = -2130
$ws.Range("H92").Value = 1196.8889
$ws.Range("I92").Value = 1490.4286
$ws.Range("K92").Value = 1490.4286
$ws.Range("M92").Value = -242.4286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1909.8
$ws.Range("I20").Value = 1717
$ws.Range("K20").Value = 1717
$ws.Range("M20").Value = -1470
$ws.Range("H33").Value = 22333.334
$ws.Range("J33").Value = 26000
$ws.Range("L33").Value = 26000
$ws.Range("N33").Value = -26672
$ws.Range("H64").Value = 696.5
$ws.Range("J64").Value = 572.25
$ws.Range("L64").Value = 572.25
$ws.Range("N64").Value = -1022.25
$ws.Range("H67").Value = 696.5
$ws.Range("J67").Value = 572.25
$ws.Range("L67").Value = 572.25
$ws.Range("N67").Value = -2132.25
$ws.Range("H99").Value = 1481.5385
$ws.Range("I99").Value = 1451.8182
$ws.Range("K99").Value = 1451.8182
$ws.Range("M99").Value = 46.18180000000007
$ws.Range("H105").Value = 2757.3333
$ws.Range("I105").Value = 2788.8
$ws.Range("K105").Value = 2788.8
$ws.Range("M105").Value = -1041.8
$ws.Range("H134").Value = 2215.25
$ws.Range("I134").Value = 2164.923
$ws.Range("K134").Value = 6494.768999999999
$ws.Range("M134").Value = -3959.768999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 11344.5
$ws.Range("I12").Value = 13461
$ws.Range("K12").Value = 13461
$ws.Range("M12").Value = -13291
$ws.Range("H22").Value = 366
$ws.Range("I22").Value = 366
$ws.Range("K22").Value = 366
$ws.Range("M22").Value = -16
$ws.Range("H31").Value = 2247.8
$ws.Range("J31").Value = 5555
$ws.Range("L31").Value = 5555
$ws.Range("N31").Value = -6145
$ws.Range("H34").Value = 2247.8
$ws.Range("J34").Value = 5555
$ws.Range("L34").Value = 5555
$ws.Range("N34").Value = -5959
$ws.Range("H35").Value = 1683.3334
$ws.Range("I35").Value = 1683.3334
$ws.Range("K35").Value = 1683.3334
$ws.Range("M35").Value = -1389.3334
$ws.Range("H62").Value = 3194.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3194.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3194.5
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -4442.5
$ws.Range("H65").Value = 3194.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3194.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 15972.5
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -22212.5
$ws.Range("H99").Value = 3666.6
$ws.Range("J99").Value = 6000
$ws.Range("L99").Value = 6000
$ws.Range("N99").Value = -8996
$ws.Range("H126").Value = 3666.6
$ws.Range("J126").Value = 6000
$ws.Range("L126").Value = 18000
$ws.Range("N126").Value = -22940
$ws.Range("H132").Value = 3168.8333
$ws.Range("I132").Value = 3189.4666
$ws.Range("K132").Value = 9568.399800000001
$ws.Range("M132").Value = -7038.399800000001
$ws.Range("H134").Value = 2290.4075
$ws.Range("I134").Value = 2312.9092
$ws.Range("J134").Value = 2191.4
$ws.Range("K134").Value = 6938.7276
$ws.Range("L134").Value = 6574.200000000001
$ws.Range("M134").Value = -4403.7276
$ws.Range("N134").Value = -11644.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 207.5
$ws.Range("J12").Value = 161.6
$ws.Range("L12").Value = 484.8
$ws.Range("N12").Value = -830.8
$ws.Range("H51").Value = 650.75
$ws.Range("I51").Value = 534.6667
$ws.Range("J51").Value = 999
$ws.Range("K51").Value = 1604.0001
$ws.Range("L51").Value = 2997
$ws.Range("M51").Value = -1144.0001
$ws.Range("N51").Value = -3917
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()
$ws.Range("H113").Value = 1507.909
$ws.Range("I113").Value = 1236
$ws.Range("K113").Value = 3708
$ws.Range("M113").Value = -1538
$ws.Range("H137").Value = 1657.8334
$ws.Range("J137").Value = 2249.3333
$ws.Range("L137").Value = 6747.999899999999
$ws.Range("N137").Value = -16947.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 90077
$ws.Range("I62").Value = 90077
$ws.Range("K62").Value = 90077
$ws.Range("M62").Value = -89391
$ws.Range("H65").Value = 90077
$ws.Range("I65").Value = 90077
$ws.Range("K65").Value = 270231
$ws.Range("M65").Value = -266799
$ws.Range("H70").Value = 5304
$ws.Range("I70").Value = 5103.5
$ws.Range("K70").Value = 5103.5
$ws.Range("M70").Value = -4833.5
$ws.Range("H73").Value = 5304
$ws.Range("I73").Value = 5103.5
$ws.Range("K73").Value = 5103.5
$ws.Range("M73").Value = -4167.5
$ws.Range("H80").Value = 2904.2856
$ws.Range("I80").Value = 2499.5
$ws.Range("J80").Value = 3066.2
$ws.Range("K80").Value = 2499.5
$ws.Range("L80").Value = 3066.2
$ws.Range("M80").Value = -1501.5
$ws.Range("N80").Value = -5062.2
$ws.Range("H83").Value = 2904.2856
$ws.Range("I83").Value = 2499.5
$ws.Range("J83").Value = 3066.2
$ws.Range("K83").Value = 12497.5
$ws.Range("L83").Value = 15331
$ws.Range("M83").Value = -7505.5
$ws.Range("N83").Value = -25315
$ws.Range("H102").Value = 2680.889
$ws.Range("I102").Value = 2812
$ws.Range("J102").Value = 2222
$ws.Range("K102").Value = 2812
$ws.Range("L102").Value = 2222
$ws.Range("M102").Value = -1190
$ws.Range("N102").Value = -5466
$ws.Range("H126").Value = 2199.4
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
$ws.Range("H132").Value = 3371.5833
$ws.Range("I132").Value = 3371.5833
$ws.Range("K132").Value = 10114.7499
$ws.Range("M132").Value = -7584.749899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2595.5
$ws.Range("J46").Value = 4279.8
$ws.Range("L46").Value = 4279.8
$ws.Range("N46").Value = -4655.8
$ws.Range("H55").Value = 926.4545000000001
$ws.Range("I55").Value = 584.2857
$ws.Range("J55").Value = 1525.25
$ws.Range("K55").Value = 584.2857
$ws.Range("L55").Value = 1525.25
$ws.Range("M55").Value = -411.2857
$ws.Range("N55").Value = -1871.25
$ws.Range("H61").Value = 958
$ws.Range("I61").Value = 919.6
$ws.Range("K61").Value = 919.6
$ws.Range("M61").Value = -717.6
$ws.Range("H74").Value = 83598.5
$ws.Range("I74").Value = 83598.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 83598.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -82600.5
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 83598.5
$ws.Range("I77").Value = 83598.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 250795.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -245803.5
$ws.Range("N77").ClearContents()
$ws.Range("H113").Value = 958
$ws.Range("I113").Value = 919.6
$ws.Range("K113").Value = 919.6
$ws.Range("M113").Value = 1250.4
$ws.Range("H136").Value = 5332.6665
$ws.Range("J136").Value = 6999
$ws.Range("L136").Value = 20997
$ws.Range("N136").Value = -26097

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 74302.336
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 74302.336
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -76174.336
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 74302.336
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 222907.008
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -232267.008
$ws.Range("H107").Value = 901.6
$ws.Range("I107").Value = 1500
$ws.Range("J107").Value = 752
$ws.Range("K107").Value = 4500
$ws.Range("L107").Value = 2256
$ws.Range("M107").Value = -2580
$ws.Range("N107").Value = -6096
$ws.Range("H126").Value = 4308.1763
$ws.Range("I126").Value = 4236.25
$ws.Range("K126").Value = 12708.75
$ws.Range("M126").Value = -10238.75
$ws.Range("H132").Value = 1691.2727
$ws.Range("I132").Value = 1776.1666
$ws.Range("J132").Value = 1309.25
$ws.Range("K132").Value = 5328.4998
$ws.Range("L132").Value = 3927.75
$ws.Range("M132").Value = -2798.4998
$ws.Range("N132").Value = -8987.75
$ws.Range("H136").Value = 2799.1428
$ws.Range("I136").Value = 3099
$ws.Range("J136").Value = 1699.6666
$ws.Range("K136").Value = 9297
$ws.Range("L136").Value = 5098.9998
$ws.Range("M136").Value = -6747
$ws.Range("N136").Value = -10198.9998
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

Write-Host "Applied all cell updates"